$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Ephb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05800433333333333
$ws.Range("H2").Value = 0.174013
$ws.Range("I2").Value = 0.02087975181349295
$ws.Range("J2").Value = 0.02087975181349295
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.306600666666667
$ws.Range("N2").Value = 3.919802
$ws.Range("O2").Value = 0.8137131711319011
$ws.Range("P2").Value = 0.8137131711319011
$ws.Range("Q2").Value = 0.07578850060288887
$ws.Range("R2").Value = 0.682096505426
$ws.Range("S2").Value = 0.01699012906060441
$ws.Range("T2").Value = 0.01699012906060441

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Ephb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05800433333333333
$ws.Range("H3").Value = 0.174013
$ws.Range("I3").Value = 0.02087975181349295
$ws.Range("J3").Value = 0.02087975181349295
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.01102233333333333
$ws.Range("N3").Value = 0.033067
$ws.Range("O3").Value = 0.006864390964089149
$ws.Range("P3").Value = 0.006864390964089149
$ws.Range("Q3").Value = 0.0006393430967777777
$ws.Range("R3").Value = 0.005754087871
$ws.Range("S3").Value = 0.000143326779680965
$ws.Range("T3").Value = 0.000143326779680965

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05800433333333333
$ws.Range("H4").Value = 0.174013
$ws.Range("I4").Value = 0.02087975181349295
$ws.Range("J4").Value = 0.02087975181349295
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2881033333333333
$ws.Range("N4").Value = 0.86431
$ws.Range("O4").Value = 0.1794224379040098
$ws.Range("P4").Value = 0.1794224379040098
$ws.Range("Q4").Value = 0.01671124178111111
$ws.Range("R4").Value = 0.15040117603
$ws.Range("S4").Value = 0.003746295973207575
$ws.Range("T4").Value = 0.003746295973207576

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Ephb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.666083666666667
$ws.Range("H5").Value = 4.998251
$ws.Range("I5").Value = 0.5997381826733804
$ws.Range("J5").Value = 0.5997381826733805
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.306600666666667
$ws.Range("N5").Value = 3.919802
$ws.Range("O5").Value = 0.8137131711319011
$ws.Range("P5").Value = 0.8137131711319011
$ws.Range("Q5").Value = 2.176906029589111
$ws.Range("R5").Value = 19.592154266302
$ws.Range("S5").Value = 0.4880148584720398
$ws.Range("T5").Value = 0.4880148584720398

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Ephb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.666083666666667
$ws.Range("H6").Value = 4.998251
$ws.Range("I6").Value = 0.5997381826733804
$ws.Range("J6").Value = 0.5997381826733805
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.01102233333333333
$ws.Range("N6").Value = 0.033067
$ws.Range("O6").Value = 0.006864390964089149
$ws.Range("P6").Value = 0.006864390964089149
$ws.Range("Q6").Value = 0.01836412953522222
$ws.Range("R6").Value = 0.165277165817
$ws.Range("S6").Value = 0.0041168373619624
$ws.Range("T6").Value = 0.004116837361962401

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Ephb1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.666083666666667
$ws.Range("H7").Value = 4.998251
$ws.Range("I7").Value = 0.5997381826733804
$ws.Range("J7").Value = 0.5997381826733805
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2881033333333333
$ws.Range("N7").Value = 0.86431
$ws.Range("O7").Value = 0.1794224379040098
$ws.Range("P7").Value = 0.1794224379040098
$ws.Range("Q7").Value = 0.4800042579788888
$ws.Range("R7").Value = 4.32003832181
$ws.Range("S7").Value = 0.1076064868393783
$ws.Range("T7").Value = 0.1076064868393783

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Efna5"
$ws.Range("C8").Value = "Ephb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.053930333333333
$ws.Range("H8").Value = 3.161791
$ws.Range("I8").Value = 0.3793820655131266
$ws.Range("J8").Value = 0.3793820655131266
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.306600666666667
$ws.Range("N8").Value = 3.919802
$ws.Range("O8").Value = 0.8137131711319011
$ws.Range("P8").Value = 0.8137131711319011
$ws.Range("Q8").Value = 1.377066076153556
$ws.Range("R8").Value = 12.393594685382
$ws.Range("S8").Value = 0.3087081835992569
$ws.Range("T8").Value = 0.3087081835992569

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Efna5"
$ws.Range("C9").Value = "Ephb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.053930333333333
$ws.Range("H9").Value = 3.161791
$ws.Range("I9").Value = 0.3793820655131266
$ws.Range("J9").Value = 0.3793820655131266
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.01102233333333333
$ws.Range("N9").Value = 0.033067
$ws.Range("O9").Value = 0.006864390964089149
$ws.Range("P9").Value = 0.006864390964089149
$ws.Range("Q9").Value = 0.01161677144411111
$ws.Range("R9").Value = 0.104550942997
$ws.Range("S9").Value = 0.002604226822445784
$ws.Range("T9").Value = 0.002604226822445784

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Efna5"
$ws.Range("C10").Value = "Ephb1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.053930333333333
$ws.Range("H10").Value = 3.161791
$ws.Range("I10").Value = 0.3793820655131266
$ws.Range("J10").Value = 0.3793820655131266
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2881033333333333
$ws.Range("N10").Value = 0.86431
$ws.Range("O10").Value = 0.1794224379040098
$ws.Range("P10").Value = 0.1794224379040098
$ws.Range("Q10").Value = 0.3036408421344445
$ws.Range("R10").Value = 2.73276757921
$ws.Range("S10").Value = 0.06806965509142393
$ws.Range("T10").Value = 0.06806965509142394
